$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving a purely-numeric-looking string are temporarily switched
# to text format before the assignment (otherwise Excel auto-converts them
# to numbers, e.g. dropping a trailing zero: "1.330" -> 1.33), then the
# explicit formatting is cleared again so the cell keeps its original
# default styling.

$ws.Range('D2').Value = '26.391.92'
$ws.Range('E2').Value = '  +0.55%  '
$ws.Range('D3').Value = '1.691.74'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('E4').Value = '  +0.44%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.97'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5467'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.94%  '
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2741'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.63%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06469'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.05'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07679'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +2.77%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.702.95'
$ws.Range('E12').Value = '  +0.79%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.548'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5837'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008403'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.75%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.35'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.03%  '
$ws.Range('D17').Value = '26.456.82'
$ws.Range('E17').Value = '  +0.57%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.951'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.41%  '
$ws.Range('E19').Value = '  +0.44%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.98'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.46'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.263'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.011'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '149.58'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.95%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1323'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +6.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.901'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.80%  '
$ws.Range('E27').Value = '  -0.64%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06352'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -4.68%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.405'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +3.82%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.330'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.610'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.36%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.597'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.685'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.99%  '
$ws.Range('E34').Value = '  +1.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6171'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.70%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.410'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.707'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.267'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.35%  '
$ws.Range('D39').Value = '1.123.45'
$ws.Range('E39').Value = '  +1.80%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01636'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8811'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.86'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.90%  '
$ws.Range('D44').Value = '1.840.17'
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000111'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -4.42%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.55'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.24%  '
$ws.Range('E47').Value = '  +0.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.236'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.77%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05280'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.141'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.71%  '
$ws.Range('E51').Value = '  +0.29%  '
